$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-11 (old extra entries), leaving header + rows 2-4
$ws.Range("A5:E11").EntireRow.Delete()

# Force the Date text column to stay plain text (not auto-converted to a
# date serial number) for the rows we are about to rewrite. The Time
# column values (e.g. "14:05:55") are not auto-recognized as dates/times
# by the engine, so they don't need this treatment.
$ws.Range("B2:B4").NumberFormat = "@"

# Update row 2
$ws.Range("A2").Value = "KL51J6070"
$ws.Range("B2").Value = "2024-11-18"
$ws.Range("C2").Value = "14:05:55"
$ws.Range("D2").Value = "IN"
$ws.Range("E2").Value = 1731918955392

# Update row 3
$ws.Range("A3").Value = "MH01s1513"
$ws.Range("B3").Value = "2024-11-18"
$ws.Range("C3").Value = "14:09:10"
$ws.Range("D3").Value = "IN"
$ws.Range("E3").Value = 1731919150890

# Update row 4
$ws.Range("A4").Value = "Mh01s1513"
$ws.Range("B4").Value = "2024-11-18"
$ws.Range("C4").Value = "14:09:34"
$ws.Range("D4").Value = "OUT"
$ws.Range("E4").Value = 1731919174136
